$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D and E to Text format so numeric-looking strings
# (e.g. "1.00", "41.639.28") are preserved exactly as text, matching
# the inline-string cell type used in the source workbook.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "41.639.28"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "2.228.27"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "251.75"
$ws.Range("E5").Value = "  +8.37%  "
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("E7").Value = "  +1.71%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "0.570"
$ws.Range("E9").Value = "  +2.21%  "
$ws.Range("E10").Value = "  +20.49%  "
$ws.Range("D11").Value = "0.0963"
$ws.Range("E11").Value = "  -2.36%  "
$ws.Range("D12").Value = "58.90"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D14").Value = "7.04"
$ws.Range("E14").Value = "  +3.60%  "
$ws.Range("D15").Value = "2.561.13"
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("D16").Value = "14.92"
$ws.Range("E16").Value = "  -1.11%  "
$ws.Range("D17").Value = "0.854"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("D18").Value = "2.219.86"
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("D19").Value = "41.604.05"
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("D22").Value = "73.04"
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("E23").Value = "  +9.94%  "
$ws.Range("D24").Value = "234.81"
$ws.Range("E24").Value = "  -1.12%  "
$ws.Range("D25").Value = "3.89"
$ws.Range("E25").Value = "  +6.69%  "
$ws.Range("E27").Value = "  +6.15%  "
$ws.Range("D28").Value = "10.43"
$ws.Range("E28").Value = "  +4.08%  "
$ws.Range("E29").Value = "  +1.67%  "
$ws.Range("D30").Value = "171.77"
$ws.Range("E30").Value = "  +1.74%  "
$ws.Range("D31").Value = "20.61"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  +1.69%  "
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("D34").Value = "5.56"
$ws.Range("E34").Value = "  +2.05%  "
$ws.Range("D35").Value = "0.0722"
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("D36").Value = "26.77"
$ws.Range("E36").Value = "  +21.74%  "
$ws.Range("D37").Value = "4.68"
$ws.Range("E37").Value = "  -1.95%  "
$ws.Range("D38").Value = "4.01"
$ws.Range("E38").Value = "  +11.64%  "
$ws.Range("D39").Value = "0.0293"
$ws.Range("E39").Value = "  +9.90%  "
$ws.Range("E40").Value = "  +2.28%  "
$ws.Range("E41").Value = "  +3.54%  "
$ws.Range("E42").Value = "  -0.41%  "
$ws.Range("D43").Value = "12.35"
$ws.Range("E43").Value = "  +23.87%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "0.209"
$ws.Range("E44").Value = "  +10.23%  "
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").Value = "5.07"
$ws.Range("E45").Value = "  +2.86%  "
$ws.Range("E46").Value = "  -2.65%  "
$ws.Range("D47").Value = "4.80"
$ws.Range("E47").Value = "  +10.72%  "
$ws.Range("E48").Value = "  +1.41%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  +7.44%  "
$ws.Range("D51").Value = "1.19"
$ws.Range("E51").Value = "  +2.03%  "
